$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date column stays stored as text (matches original "str" cell type)
# instead of being auto-converted into a date serial number by Excel.
$ws.Range("A2:A15").NumberFormat = "@"

$updates = @(
    @{ Row = 2;  Date = "2025-12-07"; Employees = "Anton Iosifov, David Cohen" }
    @{ Row = 3;  Date = "2025-12-07"; Employees = "Marina Levi, Alex Brown" }
    @{ Row = 4;  Date = "2025-12-08"; Employees = "Tommy Gun, Tony Saprano" }
    @{ Row = 5;  Date = "2025-12-08"; Employees = "Anton Iosifov, Marina Levi" }
    @{ Row = 6;  Date = "2025-12-09"; Employees = "David Cohen, Alex Brown" }
    @{ Row = 7;  Date = "2025-12-09"; Employees = "Tommy Gun, Tony Saprano" }
    @{ Row = 8;  Date = "2025-12-10"; Employees = "Anton Iosifov, Marina Levi" }
    @{ Row = 9;  Date = "2025-12-10"; Employees = "David Cohen, Alex Brown" }
    @{ Row = 10; Date = "2025-12-11"; Employees = "Tommy Gun, Tony Saprano" }
    @{ Row = 11; Date = "2025-12-11"; Employees = "Anton Iosifov, Marina Levi" }
    @{ Row = 12; Date = "2025-12-12"; Employees = "David Cohen, Alex Brown" }
    @{ Row = 13; Date = "2025-12-12"; Employees = "Tommy Gun, Tony Saprano" }
    @{ Row = 14; Date = "2025-12-13"; Employees = "Anton Iosifov, Marina Levi" }
    @{ Row = 15; Date = "2025-12-13"; Employees = "David Cohen, Alex Brown" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Date
    $ws.Cells.Item($u.Row, 3).Value = $u.Employees
}
